# AUTOMATE-7: Add patient functionality

$wb = $excel.ActiveWorkbook

# --- Guardian sheet: update the current selection ---
$guardian = $wb.Worksheets.Item("Guardian")
$guardian.Activate()
$guardian.Range("C2").Select()

# --- Patient sheet: add Date Of Birth + Mother Identifier values ---
$patient = $wb.Worksheets.Item("Patient")
$patient.Activate()

# Date Of Birth for BABY A (row 2)
$patient.Range("D2").Value = 35848
$patient.Range("D2").NumberFormat = "mm-dd-yy"

# Mother Identifier values (previously placeholder numbers)
$patient.Range("I2").Value = "MOT000001"
$patient.Range("I3").Value = "MOT000002"

# Update current selection on the Patient sheet (it is the active tab)
$patient.Range("I4").Select()
